$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current header row (row 1), pushing
# everything (header + all data rows) down by two rows.
$ws.Rows("1:2").Insert()

# --- New row 1: numeric column-index header (0..12), same bold/centered/
#     bordered style the original text header (now row 3) used to carry. ---
$headerNums = @(0,1,2,3,4,5,6,7,8,9,10,11,12)
for ($i = 0; $i -lt 13; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headerNums[$i]
}
$hdrRange = $ws.Range("A1:M1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108   # xlCenter
$hdrRange.VerticalAlignment = -4160     # xlTop
$hdrRange.Borders.LineStyle = 1         # xlContinuous

# --- New row 2: mostly blank "Flange" / "Drive" sub-header row. ---
$ws.Range("A2:M2").Value = ""
$ws.Cells.Item(2, 3).Value = "Flange"   # C2
$ws.Cells.Item(2, 6).Value = "Drive"    # F2

# --- Row 3 (the original text header, shifted down by the insert) loses
#     the bold/border/center formatting it used to have. ---
$ws.Range("A3:M3").ClearFormats()

$wb.Save()
